# Adecuacion de la implementacion de lectura excel para credito empresarial
#
# Inserts a "Moneda" column (currency symbol) right before the existing
# "Monto" column, and a "Comentarios Ratios" column right before
# "Observaciones aprobación". Also converts a handful of numeric-looking
# columns (expediente id, Dias Pago, Tasa Preferencial, Monto, count) to
# text-formatted cells (NumberFormat "@"), matching the updated reader
# expectations.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Insert the two new columns.
#    Before: A..S = Cliente..Contraseña (19 cols), " Monto" in G.
#    Insert a blank column before G -> new G is blank ("Moneda"),
#    old G.."S" shift to H.."T".
# ---------------------------------------------------------------------
$ws.Columns("G").Insert()

# After the shift above, "Observaciones aprobación" (old R) now sits in S.
# Insert another blank column right before it for "Comentarios Ratios".
$ws.Columns("S").Insert()

try { $ws.Columns("G").ColumnWidth = 8.6 } catch {}
try { $ws.Columns("S").ColumnWidth = 29.14 } catch {}

# ---------------------------------------------------------------------
# 2. New "Moneda" column (currency symbol), data rows first then header -
#    mirrors the order new labels were introduced in the source edit.
# ---------------------------------------------------------------------
foreach ($r in 2, 3) {
    # Formatted like its neighbours (bold, vertically centred) by copying
    # the format from column F.
    $ws.Range("F$r").Copy() | Out-Null
    $ws.Range("G$r").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
    $ws.Range("G$r").Value = "S/"
}
$excel.CutCopyMode = $false

$ws.Range("G1").Value = "Moneda"

# Expediente id (col A): was a plain number, now stored as text.
$ws.Range("A2:A3").NumberFormat = "@"
$ws.Range("A2").Value = "2363"
$ws.Range("A3").Value = "2517"

# Dias Pago / Tasa Preferencial (cols L, M after the insert): were plain
# numbers, now stored as text (still bold + vertically centred).
$ws.Range("L2:M3").NumberFormat = "@"
$ws.Range("L2").Value = "90"
$ws.Range("L3").Value = "90"
$ws.Range("M2").Value = "30"
$ws.Range("M3").Value = "30"

# New "Comentarios Ratios" column: same text as "Nota de operación"
# (column R), left with default (unstyled) formatting.
$ws.Range("S1").Value = "Comentarios Ratios"
foreach ($r in 2, 3) {
    $ws.Range("S$r").Value = $ws.Range("R$r").Value2
}

# Monto / count (cols H, I after the insert): keep the numeric values but
# switch the cell format to text.
$ws.Range("H2:I3").NumberFormat = "@"

# ---------------------------------------------------------------------
# 4. Selection / view cosmetics (best effort).
# ---------------------------------------------------------------------
$ws.Range("S7").Select() | Out-Null
